$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
# row 40 (G40=5505)
$ws.Cells.Item(40, 8).Value = 3960
$ws.Cells.Item(40, 10).Value = 4413.3335
$ws.Cells.Item(40, 12).Value = 4413.3335
$ws.Cells.Item(40, 14).Value = -4763.3335
# row 51 (G51=5486)
$ws.Cells.Item(51, 8).Value = 3884
$ws.Cells.Item(51, 9).Value = 2986.3333
$ws.Cells.Item(51, 10).Value = 4781.6665
$ws.Cells.Item(51, 11).Value = 2986.3333
$ws.Cells.Item(51, 12).Value = 4781.6665
$ws.Cells.Item(51, 13).Value = -2502.3333
$ws.Cells.Item(51, 14).Value = -5749.6665
# row 92 (G92=19901)
$ws.Cells.Item(92, 8).Value = 554.72
$ws.Cells.Item(92, 9).Value = 577.2632
$ws.Cells.Item(92, 10).Value = 483.33334
$ws.Cells.Item(92, 11).Value = 577.2632
$ws.Cells.Item(92, 12).Value = 483.33334
$ws.Cells.Item(92, 13).Value = 670.7368
$ws.Cells.Item(92, 14).Value = -2979.33334
# row 94 (G94=19905)
$ws.Cells.Item(94, 8).Value = 5249.364
$ws.Cells.Item(94, 9).Value = 963.2857
$ws.Cells.Item(94, 10).Value = 12750
$ws.Cells.Item(94, 11).Value = 963.2857
$ws.Cells.Item(94, 12).Value = 12750
$ws.Cells.Item(94, 13).Value = -512.2857
$ws.Cells.Item(94, 14).Value = -13652
# row 96 (G96=19894)
$ws.Cells.Item(96, 8).Value = 726.9
$ws.Cells.Item(96, 9).Value = 586
$ws.Cells.Item(96, 10).Value = 1995
$ws.Cells.Item(96, 11).Value = 1758
$ws.Cells.Item(96, 12).Value = 5985
$ws.Cells.Item(96, 13).Value = -385
$ws.Cells.Item(96, 14).Value = -8731
# row 97 (G97=19885)
$ws.Cells.Item(97, 8).Value = 25002224
$ws.Cells.Item(97, 10).Value = 25002224
$ws.Cells.Item(97, 12).Value = 75006672
$ws.Cells.Item(97, 14).Value = -75007664
# row 125 (G125=36228)
$ws.Cells.Item(125, 8).Value = 111112790
$ws.Cells.Item(125, 9).Value = 166668000
$ws.Cells.Item(125, 10).Value = 2366.3333
$ws.Cells.Item(125, 11).Value = 1500012000
$ws.Cells.Item(125, 12).Value = 21296.9997
$ws.Cells.Item(125, 13).Value = -1500009540
$ws.Cells.Item(125, 14).Value = -26216.9997
# row 129 (G129=36115)
$ws.Cells.Item(129, 8).Value = 1546.7727
$ws.Cells.Item(129, 9).Value = 1083.2727
$ws.Cells.Item(129, 11).Value = 3249.8181
$ws.Cells.Item(129, 13).Value = 1750.1819
# row 137 (G137=44013)
$ws.Cells.Item(137, 8).Value = 4018.7021
$ws.Cells.Item(137, 9).Value = 2010.0333
$ws.Cells.Item(137, 10).Value = 7563.4116
$ws.Cells.Item(137, 11).Value = 6030.0999
$ws.Cells.Item(137, 12).Value = 22690.2348
$ws.Cells.Item(137, 13).Value = -3480.0999
$ws.Cells.Item(137, 14).Value = -27790.2348
# row 138 (G138=44169)
$ws.Cells.Item(138, 8).Value = 1758857.4
$ws.Cells.Item(138, 9).Value = 2900.1177
$ws.Cells.Item(138, 10).Value = 2505139
$ws.Cells.Item(138, 11).Value = 8700.3531
$ws.Cells.Item(138, 12).Value = 7515417
$ws.Cells.Item(138, 13).Value = -3560.3531
$ws.Cells.Item(138, 14).Value = -7525697

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
# row 32 (G32=44147)
$ws.Cells.Item(32, 8).Value = 2903943.8
$ws.Cells.Item(32, 9).Value = 3128495
$ws.Cells.Item(32, 11).Value = 3128495
$ws.Cells.Item(32, 13).Value = -3128208
# row 102 (G102=19945)
$ws.Cells.Item(102, 8).Value = 26320148
$ws.Cells.Item(102, 9).Value = 45458184
$ws.Cells.Item(102, 10).Value = 5350
$ws.Cells.Item(102, 11).Value = 45458184
$ws.Cells.Item(102, 12).Value = 5350
$ws.Cells.Item(102, 13).Value = -45456562
$ws.Cells.Item(102, 14).Value = -8594

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
# row 99 (G99=19943)
$ws.Cells.Item(99, 8).Value = 13496221
$ws.Cells.Item(99, 9).Value = 13890210
$ws.Cells.Item(99, 10).Value = 12989665
$ws.Cells.Item(99, 11).Value = 13890210
$ws.Cells.Item(99, 12).Value = 12989665
$ws.Cells.Item(99, 13).Value = -13888712
$ws.Cells.Item(99, 14).Value = -12992661
# row 105 (G105=19947)
$ws.Cells.Item(105, 8).Value = 3594.5144
$ws.Cells.Item(105, 9).Value = 2773.7827
$ws.Cells.Item(105, 11).Value = 2773.7827
$ws.Cells.Item(105, 13).Value = -1026.7827
# row 117 (G117=26124)
$ws.Cells.Item(117, 8).Value = 95899.5
$ws.Cells.Item(117, 10).Value = 95899.5
$ws.Cells.Item(117, 12).Value = 95899.5
$ws.Cells.Item(117, 14).Value = -105077.5
# row 134 (G134=43998)
$ws.Cells.Item(134, 8).Value = 7582576.5
$ws.Cells.Item(134, 10).Value = 11043.0625
$ws.Cells.Item(134, 12).Value = 33129.1875
$ws.Cells.Item(134, 14).Value = -38199.1875

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
# row 31 (G31=44023)
$ws.Cells.Item(31, 8).Value = 7131.744
$ws.Cells.Item(31, 9).Value = 2067.842
$ws.Cells.Item(31, 10).Value = 11140.667
$ws.Cells.Item(31, 11).Value = 2067.842
$ws.Cells.Item(31, 12).Value = 11140.667
$ws.Cells.Item(31, 13).Value = -1772.842
$ws.Cells.Item(31, 14).Value = -11730.667
# row 34 (G34=44023)
$ws.Cells.Item(34, 8).Value = 7131.744
$ws.Cells.Item(34, 9).Value = 2067.842
$ws.Cells.Item(34, 10).Value = 11140.667
$ws.Cells.Item(34, 11).Value = 2067.842
$ws.Cells.Item(34, 12).Value = 11140.667
$ws.Cells.Item(34, 13).Value = -1865.842
$ws.Cells.Item(34, 14).Value = -11544.667
# row 62 (G62=12580)
$ws.Cells.Item(62, 8).Value = 18524040
$ws.Cells.Item(62, 9).Value = 33337872
$ws.Cells.Item(62, 10).Value = 6748.5
$ws.Cells.Item(62, 11).Value = 33337872
$ws.Cells.Item(62, 12).Value = 6748.5
$ws.Cells.Item(62, 13).Value = -33337248
$ws.Cells.Item(62, 14).Value = -7996.5
# row 65 (G65=12580)
$ws.Cells.Item(65, 8).Value = 18524040
$ws.Cells.Item(65, 9).Value = 33337872
$ws.Cells.Item(65, 10).Value = 6748.5
$ws.Cells.Item(65, 11).Value = 166689360
$ws.Cells.Item(65, 12).Value = 33742.5
$ws.Cells.Item(65, 13).Value = -166686240
$ws.Cells.Item(65, 14).Value = -39982.5

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
# row 60 (G60=4750)
$ws.Cells.Item(60, 8).Value = 495
$ws.Cells.Item(60, 9).Value = 495
$ws.Cells.Item(60, 11).Value = 1485
$ws.Cells.Item(60, 13).Value = -1234
# row 95 (G95=19838)
$ws.Cells.Item(95, 8).Value = 4262
$ws.Cells.Item(95, 10).Value = 4500
$ws.Cells.Item(95, 12).Value = 13500
$ws.Cells.Item(95, 14).Value = -17618
# row 107 (G107=27838)
$ws.Cells.Item(107, 8).Value = 12223026
$ws.Cells.Item(107, 9).Value = 2500511.8
$ws.Cells.Item(107, 11).Value = 7501535.399999999
$ws.Cells.Item(107, 13).Value = -7499615.399999999
# row 122 (G122=36078)
$ws.Cells.Item(122, 8).Value = 4041741.5
$ws.Cells.Item(122, 10).Value = 1004
$ws.Cells.Item(122, 12).Value = 9036
$ws.Cells.Item(122, 14).Value = -13936

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
# row 80 (G80=12521)
$ws.Cells.Item(80, 8).Value = 2301.9092
$ws.Cells.Item(80, 9).Value = 2391.5557
$ws.Cells.Item(80, 10).Value = 1898.5
$ws.Cells.Item(80, 11).Value = 2391.5557
$ws.Cells.Item(80, 12).Value = 1898.5
$ws.Cells.Item(80, 13).Value = -1393.5557
$ws.Cells.Item(80, 14).Value = -3894.5
# row 83 (G83=12521)
$ws.Cells.Item(83, 8).Value = 2301.9092
$ws.Cells.Item(83, 9).Value = 2391.5557
$ws.Cells.Item(83, 10).Value = 1898.5
$ws.Cells.Item(83, 11).Value = 11957.7785
$ws.Cells.Item(83, 12).Value = 9492.5
$ws.Cells.Item(83, 13).Value = -6965.7785
$ws.Cells.Item(83, 14).Value = -19476.5
# row 92 (G92=18094)
$ws.Cells.Item(92, 8).Value = 19300
$ws.Cells.Item(92, 10).Value = 19300
$ws.Cells.Item(92, 12).Value = 19300
$ws.Cells.Item(92, 14).Value = -23044
# row 113 (G113=27710)
$ws.Cells.Item(113, 8).Value = 6468.5674
$ws.Cells.Item(113, 9).Value = 4251.4116
$ws.Cells.Item(113, 11).Value = 4251.4116
$ws.Cells.Item(113, 13).Value = -2081.4116
# row 122 (G122=36182)
$ws.Cells.Item(122, 8).Value = 7266056.5
$ws.Cells.Item(122, 9).Value = 9081821
$ws.Cells.Item(122, 11).Value = 27245463
$ws.Cells.Item(122, 13).Value = -27243013
# row 126 (G126=36184)
$ws.Cells.Item(126, 8).Value = 4272.6
$ws.Cells.Item(126, 9).Value = 2513.6667
$ws.Cells.Item(126, 10).Value = 6911
$ws.Cells.Item(126, 11).Value = 7541.000100000001
$ws.Cells.Item(126, 12).Value = 20733
$ws.Cells.Item(126, 13).Value = -5071.000100000001
$ws.Cells.Item(126, 14).Value = -25673
# row 132 (G132=44008)
$ws.Cells.Item(132, 8).Value = 4694.731
$ws.Cells.Item(132, 9).Value = 1782.8889
$ws.Cells.Item(132, 11).Value = 5348.6667
$ws.Cells.Item(132, 13).Value = -2818.6667
# row 133 (G133=41854)
$ws.Cells.Item(133, 8).Value = 60564
$ws.Cells.Item(133, 10).Value = 60564
$ws.Cells.Item(133, 12).Value = 60564
$ws.Cells.Item(133, 14).Value = -70684

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
# row 7 (G7=36249)
$ws.Cells.Item(7, 8).Value = 4359.684
$ws.Cells.Item(7, 9).Value = 3914.4443
$ws.Cells.Item(7, 11).Value = 3914.4443
$ws.Cells.Item(7, 13).Value = -3802.4443
# row 40 (G40=36248)
$ws.Cells.Item(40, 8).Value = 6019.706
$ws.Cells.Item(40, 9).Value = 2764.5715
$ws.Cells.Item(40, 11).Value = 2764.5715
$ws.Cells.Item(40, 13).Value = -2628.5715
# row 46 (G46=5282)
$ws.Cells.Item(46, 8).Value = 2432.2812
$ws.Cells.Item(46, 10).Value = 3182.4119
$ws.Cells.Item(46, 12).Value = 3182.4119
$ws.Cells.Item(46, 14).Value = -3558.4119
# row 82 (G82=12565)
$ws.Cells.Item(82, 8).Value = 626510.3
$ws.Cells.Item(82, 10).Value = 2284.875
$ws.Cells.Item(82, 12).Value = 2284.875
$ws.Cells.Item(82, 14).Value = -3006.875
# row 85 (G85=12565)
$ws.Cells.Item(85, 8).Value = 626510.3
$ws.Cells.Item(85, 10).Value = 2284.875
$ws.Cells.Item(85, 12).Value = 2284.875
$ws.Cells.Item(85, 14).Value = -4780.875
# row 88 (G88=10961)
$ws.Cells.Item(88, 8).Value = 53999.5
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 53999.5
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 53999.5
$ws.Cells.Item(88, 13).ClearContents()
$ws.Cells.Item(88, 14).Value = -54855.5
# row 91 (G91=10961)
$ws.Cells.Item(91, 8).Value = 53999.5
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 53999.5
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 53999.5
$ws.Cells.Item(91, 13).ClearContents()
$ws.Cells.Item(91, 14).Value = -56963.5
# row 122 (G122=36247)
$ws.Cells.Item(122, 8).Value = 4331.4165
$ws.Cells.Item(122, 10).Value = 7187
$ws.Cells.Item(122, 12).Value = 21561
$ws.Cells.Item(122, 14).Value = -26461
# row 126 (G126=36249)
$ws.Cells.Item(126, 8).Value = 4359.684
$ws.Cells.Item(126, 9).Value = 3914.4443
$ws.Cells.Item(126, 11).Value = 11743.3329
$ws.Cells.Item(126, 13).Value = -9273.332900000001

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
# row 2 (G2=3307)
$ws.Cells.Item(2, 8).Value = 7873
$ws.Cells.Item(2, 9).Value = 7873
$ws.Cells.Item(2, 11).Value = 7873
$ws.Cells.Item(2, 13).Value = -7761
# row 4 (G4=2996)
$ws.Cells.Item(4, 8).Value = 9980.714
$ws.Cells.Item(4, 10).Value = 5067
$ws.Cells.Item(4, 12).Value = 5067
$ws.Cells.Item(4, 14).Value = -5293
# row 107 (G107=27746)
$ws.Cells.Item(107, 8).Value = 23811158
$ws.Cells.Item(107, 9).Value = 949.75
$ws.Cells.Item(107, 11).Value = 2849.25
$ws.Cells.Item(107, 13).Value = -929.25
# row 113 (G113=27752)
$ws.Cells.Item(113, 8).Value = 755.7320999999999
$ws.Cells.Item(113, 9).Value = 686.5789
$ws.Cells.Item(113, 10).Value = 901.7222
$ws.Cells.Item(113, 11).Value = 2059.7367
$ws.Cells.Item(113, 12).Value = 2705.1666
$ws.Cells.Item(113, 13).Value = 110.2633000000001
$ws.Cells.Item(113, 14).Value = -7045.1666
# row 122 (G122=36208)
$ws.Cells.Item(122, 8).Value = 178612.73
$ws.Cells.Item(122, 9).Value = 253675.5
$ws.Cells.Item(122, 10).Value = 7040.7144
$ws.Cells.Item(122, 11).Value = 761026.5
$ws.Cells.Item(122, 12).Value = 21122.1432
$ws.Cells.Item(122, 13).Value = -758576.5
$ws.Cells.Item(122, 14).Value = -26022.1432
